$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel COM auto-converts them to numeric
# values and the exact string formatting (e.g. trailing zeros) is lost.
$numericLookingCells = @(
    "D5", "D6", "D8", "D10", "D11", "D14", "D18", "D19", "D20", "D21", "D22", "D23", "D27", "D28", "D30", "D33", "D34", "D35", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50"
)
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '69.424.11'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').Value = '3.668.76'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '644.60'
$ws.Range('E5').Value = '  -5.23%  '
$ws.Range('D6').Value = '158.85'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.498'
$ws.Range('E8').Value = '  +0.73%  '
$ws.Range('E9').Value = '  -1.26%  '
$ws.Range('D10').Value = '7.10'
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('D11').Value = '0.441'
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('E12').Value = '  -1.32%  '
$ws.Range('D13').Value = '4.285.32'
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('D14').Value = '32.51'
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').Value = '3.671.09'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('D16').Value = '69.396.40'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '15.90'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').Value = '6.45'
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('D20').Value = '467.59'
$ws.Range('E20').Value = '  -1.51%  '
$ws.Range('D21').Value = '10.03'
$ws.Range('E21').Value = '  +2.24%  '
$ws.Range('D22').Value = '0.645'
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('D23').Value = '79.28'
$ws.Range('E23').Value = '  -1.15%  '
$ws.Range('D24').Value = '3.814.29'
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('E26').Value = '  -1.22%  '
$ws.Range('D27').Value = '10.87'
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('D28').Value = '8.98'
$ws.Range('E28').Value = '  -1.64%  '
$ws.Range('E29').Value = '  -3.32%  '
$ws.Range('D30').Value = '1.69'
$ws.Range('E30').Value = '  -2.96%  '
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('E32').Value = '  -1.34%  '
$ws.Range('D33').Value = '26.69'
$ws.Range('E33').Value = '  -1.16%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '6.40'
$ws.Range('E34').Value = '  -2.78%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '0.163'
$ws.Range('E35').Value = '  +0.35%  '
$ws.Range('D36').Value = '3.665.75'
$ws.Range('E36').Value = '  -0.48%  '
$ws.Range('D37').Value = '8.41'
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').Value = '5.85'
$ws.Range('E39').Value = '  -5.53%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').Value = '178.48'
$ws.Range('E40').Value = '  +5.99%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').Value = '2.19'
$ws.Range('E42').Value = '  -3.21%  '
$ws.Range('D43').Value = '0.0892'
$ws.Range('E43').Value = '  -1.41%  '
$ws.Range('D44').Value = '0.923'
$ws.Range('E44').Value = '  -1.78%  '
$ws.Range('D45').Value = '47.12'
$ws.Range('E45').Value = '  +0.59%  '
$ws.Range('D46').Value = '28.68'
$ws.Range('E46').Value = '  +1.22%  '
$ws.Range('D47').Value = '2.67'
$ws.Range('E47').Value = '  -1.95%  '
$ws.Range('E48').Value = '  -2.33%  '
$ws.Range('D49').Value = '7.78'
$ws.Range('E49').Value = '  -1.13%  '
$ws.Range('D50').Value = '0.000264'
$ws.Range('E50').Value = '  -5.54%  '
$ws.Range('E51').Value = '  -5.55%  '
